$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handoff
#
# A new source file (e1d43867-16fc-444c-9f94-ce2de3fc887a.md) was handed off
# for localization. It needs a new row inserted just above the trailing
# ".localization-config" row (row 8) on all three sheets (Overview, zh-cn,
# de-de) -- pushing ".localization-config" down to row 9 -- and populated
# with the new file's status / handoff info.
# ---------------------------------------------------------------------------

$newFile   = "e1d43867-16fc-444c-9f94-ce2de3fc887a.md"
$newFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/a1ae0ff163335b0e13d7720b52528afb358a7617/e2e/e1d43867-16fc-444c-9f94-ce2de3fc887a.md"
$cfgFile   = ".localization-config"
$cfgUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/a1ae0ff163335b0e13d7720b52528afb358a7617/.localization-config"

# ===========================================================================
# Sheet 1: Overview
# ===========================================================================
$ws = $wb.Worksheets.Item("Overview")

$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = $newFile
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = "Ready for handoff"

# Row 9 (shifted .localization-config row) already carries over the correct
# values/styles from the insert -- nothing else to set there.

# Fix up hyperlinks: A8 still carries the OLD (.localization-config) link
# after the row insert (hyperlink ranges don't shift), so re-point it at the
# newly-inserted file, then add the missing link for the shifted-down
# .localization-config row at A9.
$ws.Hyperlinks.Add($ws.Range("A8"), $newFileUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("A9"), $cfgUrl, "", "", $cfgFile)

# ===========================================================================
# Sheet 2: zh-cn
# ===========================================================================
$ws = $wb.Worksheets.Item("zh-cn")

$xlfZh    = "e1d43867-16fc-444c-9f94-ce2de3fc887a.ca15275c8f9fd0320d66b08c9bc5a73c72e4ccb6.zh-cn.xlf"
$xlfZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ea7a53ff93e00b51774037e00e0cbfc5db7e098c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/e1d43867-16fc-444c-9f94-ce2de3fc887a.ca15275c8f9fd0320d66b08c9bc5a73c72e4ccb6.zh-cn.xlf"

$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = $newFile
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = $xlfZh
$ws.Range("D8").Value = "2016-03-04 05:39:44"
$ws.Range("G8").Value = "0001-01-01 00:00:00"
$ws.Range("H8").Value = "Include"

# Row 9 (shifted .localization-config row) already carries over the correct
# values/styles from the insert -- nothing else to set there.

$ws.Hyperlinks.Add($ws.Range("A8"), $newFileUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("C8"), $xlfZhUrl, "", "", $xlfZh)
$ws.Hyperlinks.Add($ws.Range("A9"), $cfgUrl, "", "", $cfgFile)

# ===========================================================================
# Sheet 3: de-de
# ===========================================================================
$ws = $wb.Worksheets.Item("de-de")

$xlfDe    = "e1d43867-16fc-444c-9f94-ce2de3fc887a.ca15275c8f9fd0320d66b08c9bc5a73c72e4ccb6.de-de.xlf"
$xlfDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/68c0be46656e79480df10f33d62d91c49d38b07e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/e1d43867-16fc-444c-9f94-ce2de3fc887a.ca15275c8f9fd0320d66b08c9bc5a73c72e4ccb6.de-de.xlf"

$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = $newFile
$ws.Range("B8").Value = "Ready for handoff"
$ws.Range("C8").Value = $xlfDe
$ws.Range("D8").Value = "2016-03-04 05:39:58"
$ws.Range("G8").Value = "0001-01-01 00:00:00"
$ws.Range("H8").Value = "Include"

# Row 9 (shifted .localization-config row) already carries over the correct
# values/styles from the insert -- nothing else to set there.

$ws.Hyperlinks.Add($ws.Range("A8"), $newFileUrl, "", "", $newFile)
$ws.Hyperlinks.Add($ws.Range("C8"), $xlfDeUrl, "", "", $xlfDe)
$ws.Hyperlinks.Add($ws.Range("A9"), $cfgUrl, "", "", $cfgFile)
